# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data column is AC = column index 29. We never touch column A (the
# sequential row index), only columns B..AC.
$firstCol = 2
$lastCol = 29

function Get-RowValues($sheet, $row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $sheet.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($sheet, $row, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $c = $firstCol + $i
        $sheet.Cells.Item($row, $c).Value = $vals[$i]
    }
}

function Swap-Rows($sheet, $rowA, $rowB) {
    $valsA = Get-RowValues $sheet $rowA
    $valsB = Get-RowValues $sheet $rowB
    Set-RowValues $sheet $rowA $valsB
    Set-RowValues $sheet $rowB $valsA
}

function Copy-RowInto($sheet, $srcRow, $dstRow) {
    $vals = Get-RowValues $sheet $srcRow
    Set-RowValues $sheet $dstRow $vals
}

# Re-pair fixtures that were attributed to the wrong match id (rows keep
# their sequential "A" index, but the rest of the row content swaps).
Swap-Rows $ws 200 201
Swap-Rows $ws 202 203
Swap-Rows $ws 204 205
Swap-Rows $ws 206 207

# Rows 241 and 242 hold the corrected data for matches currently recorded
# in rows 237 and 238; copy it over before removing the now-duplicated
# trailing rows.
Copy-RowInto $ws 241 237
Copy-RowInto $ws 242 238

# Remove the now redundant trailing rows (237-240 duplicates / superseded
# entries), which also fixes the sheet dimension and shared-string count.
$ws.Range("A239:A242").EntireRow.Delete()
